$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlink on B2 before clearing/rewriting the sheet.
$ws.Range("B2").Hyperlinks.Delete()

# Wipe all existing cell content (old CP004..CP020 rows, Dato006..Dato00N columns, etc.)
$ws.Cells.Clear()

# --- Header row ---
$ws.Range("A1").Value = "TituloCPs"
$ws.Range("B1").Value = "Dato001"
$ws.Range("C1").Value = "Dato002"
$ws.Range("D1").Value = "Dato003"
$ws.Range("E1").Value = "Dato004"
$ws.Range("F1").Value = "Dato005"

# --- Row 2: CP001_login_fallido ---
$ws.Range("A2").Value = "CP001_login_fallido"
$ws.Range("B2").Value = "jisola.tsoft@gmail.com"
$ws.Range("B2").Style = "Hipervínculo"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:jisola.tsoft@gmail.com")
$ws.Range("C2").Value = 12345678
$ws.Range("D2").Value = "¿Olvidaste tu contraseña?"

# --- Row 3: CP002_login_exitoso ---
$ws.Range("A3").Value = "CP002_login_exitoso"
$ws.Range("B3").Value = "jisola.tsoft@gmail.com"
$ws.Range("B3").Style = "Hipervínculo"
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:jisola.tsoft@gmail.com")
$ws.Range("C3").Value = 12061990
$ws.Range("D3").Value = "Te damos la bienvenida a Facebook, Juan"

# --- Row 4: CP003_cerrar_sesion ---
$ws.Range("A4").Value = "CP003_cerrar_sesion"
$ws.Range("B4").Value = "jisola.tsoft@gmail.com"
$ws.Range("B4").Style = "Hipervínculo"
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:jisola.tsoft@gmail.com")
$ws.Range("C4").Value = 12061990

# --- Row 5: CP004_modo_oscuro ---
$ws.Range("A5").Value = "CP004_modo_oscuro"

# --- Row 6: CP005_buscar_persona ---
$ws.Range("A6").Value = "CP005_buscar_persona"

# --- Row 7: CP006_enviar_solicitud ---
$ws.Range("A7").Value = "CP006_enviar_solicitud"

# --- Row 8: CP007_cancelar_solicitud ---
$ws.Range("A8").Value = "CP007_cancelar_solicitud"

# --- Row 9: CP008_meGusta_pagina ---
$ws.Range("A9").Value = "CP008_meGusta_pagina"
$ws.Range("B9").Value = "jisola.tsoft@gmail.com"
$ws.Range("B9").Style = "Hipervínculo"
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:jisola.tsoft@gmail.com")
$ws.Range("C9").Value = 12061990
$ws.Range("D9").Value = "Tsoft"

# --- Row 10: CP009_crear_publicacion ---
$ws.Range("A10").Value = "CP009_crear_publicacion"

# --- Row 11: CP010_crear_historia ---
$ws.Range("A11").Value = "CP010_crear_historia"

# --- Row 12: CP011_enviar_mensaje ---
$ws.Range("A12").Value = "CP011_enviar_mensaje"

# Restore column widths (Cells.Clear() strips formatting, but widths are column-level;
# re-assert them to stay in sync with the target layout).
$ws.Columns.Item(1).ColumnWidth = 35.28515625
$ws.Columns.Item(2).ColumnWidth = 13.7109375
$ws.Columns.Item(4).ColumnWidth = 38.42578125

# Page orientation
$ws.PageSetup.Orientation = 1

# Selection
$null = $ws.Range("D10").Select()
